$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.994.64"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.871.34"
$ws.Range("E3").Value = "  -2.64%  "
$ws.Range("D5").Value = "'319.45"
$ws.Range("E5").Value = "  -3.68%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("D7").Value = "'0.5094"
$ws.Range("E7").Value = "  -2.97%  "
$ws.Range("D8").Value = "'0.3952"
$ws.Range("E8").Value = "  -2.72%  "
$ws.Range("D9").Value = "'0.08214"
$ws.Range("E9").Value = "  -3.77%  "
$ws.Range("D10").Value = "'42.19"
$ws.Range("E10").Value = "  -2.10%  "
$ws.Range("E11").Value = "  -3.22%  "
$ws.Range("D12").Value = "'23.54"
$ws.Range("E12").Value = "  +4.33%  "
$ws.Range("D13").Value = "1.866.14"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").Value = "'6.301"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("E15").Value = "  -2.85%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "'91.94"
$ws.Range("E17").Value = "  -4.97%  "
$ws.Range("D18").Value = "'0.00001086"
$ws.Range("E18").Value = "  -2.79%  "
$ws.Range("D19").Value = "'0.06390"
$ws.Range("E19").Value = "  -4.76%  "
$ws.Range("D20").Value = "'17.95"
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "29.980.41"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").Value = "'5.837"
$ws.Range("E23").Value = "  -3.74%  "
$ws.Range("D24").Value = "'11.13"
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("E25").Value = "  -2.27%  "
$ws.Range("D26").Value = "2.083.94"
$ws.Range("E26").Value = "  -2.65%  "
$ws.Range("D27").Value = "'160.67"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "'21.02"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("E29").Value = "  -9.31%  "
$ws.Range("D30").Value = "'127.72"
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("D31").Value = "'1.069"
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("D33").Value = "'5.938"
$ws.Range("E33").Value = "  -3.17%  "
$ws.Range("D34").Value = "'3.716"
$ws.Range("E34").Value = "  +1.82%  "
$ws.Range("D36").Value = "'5.217"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "'0.06367"
$ws.Range("E37").Value = "  -3.61%  "
$ws.Range("D38").Value = "'0.2143"
$ws.Range("E38").Value = "  -4.01%  "
$ws.Range("D39").Value = "'1.178"
$ws.Range("E39").Value = "  -4.85%  "
$ws.Range("D40").Value = "'8.549"
$ws.Range("E40").Value = "  -6.18%  "
$ws.Range("D41").Value = "'0.6312"
$ws.Range("E41").Value = "  -3.96%  "
$ws.Range("D42").Value = "'11.35"
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("D43").Value = "'1.202"
$ws.Range("E43").Value = "  -3.69%  "
$ws.Range("D44").Value = "'0.9994"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'12.97"
$ws.Range("E45").Value = "  -3.02%  "
$ws.Range("D46").Value = "'0.5910"
$ws.Range("E46").Value = "  -4.75%  "
$ws.Range("D47").Value = "'3.641"
$ws.Range("E47").Value = "  -3.88%  "
$ws.Range("D48").Value = "'2.012"
$ws.Range("E48").Value = "  -3.81%  "
$ws.Range("D49").Value = "'122.64"
$ws.Range("E49").Value = "  -2.16%  "
$ws.Range("D50").Value = "'1.205"
$ws.Range("E50").Value = "  -3.49%  "
$ws.Range("D51").Value = "'1.124"
$ws.Range("E51").Value = "  -2.81%  "
